$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 47548
$ws.Range("B3").Value = 92.43986077227223
$ws.Range("B4").Value = 10.68069421272411
$ws.Range("B5").Value = 47.44
$ws.Range("B6").Value = 89.0675
$ws.Range("B7").Value = 97.91
